$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in missing details for existing row 6 (Sayan Basak) ---
# Mobile No
$ws.Range("C6").Value = "70035 97510"

# Email (hyperlink, mailto)
$ws.Range("D6").Value = "anodiam.sb@gmail.com"
$ws.Range("D6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:anodiam.sb@gmail.com") | Out-Null

# --- New row 7: Anupam Sen ---
$ws.Range("A7").Value = 6
$ws.Range("A7").HorizontalAlignment = $ws.Range("A6").HorizontalAlignment

$ws.Range("B7").Value = "Anupam Sen"

$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "10415829829"

$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "SBIN0001486"

$ws.Range("G7").Value = "State Bank Of India"

# --- New header columns J (Branch) and K (Branch Address) ---
$ws.Range("J1").Value = "Branch"
$ws.Range("J1").Font.Bold = $ws.Range("F1").Font.Bold

$ws.Range("J7").Value = "Kasba"

$ws.Range("K1").Value = "Branch Address"
$ws.Range("K1").Font.Bold = $ws.Range("F1").Font.Bold

$ws.Range("K7").Value = "250, B B CHATTERJEE ROAD, KASBA, KOLKATA - 700042"

$ws.Range("C7").Value = "98303 04429"

$ws.Range("E7").Value = "Economics"

$ws.Range("F7").Value = "Active"

$ws.Range("J6").Value = "College Street"

# Widen column G to fit the new, longer bank name
$ws.Columns.Item(7).ColumnWidth = 15.75

# Move the active selection, as recorded in the saved workbook
$ws.Range("J7").Select() | Out-Null
